$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 35; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE

    $currentPeriod = $hCell.Value2
    if ($currentPeriod -ne $null) {
        $hCell.Value2 = $currentPeriod - 1
    }

    $iCell.Value2 = "'04-Nov-2025"
}
